$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Digits Kinematics Length Data: refresh the DP_Length (m) column
# (column A) values for rows 2-6 with the newly computed lengths.
$ws.Range("A2").Value = 0.029360877984147542
$ws.Range("A3").Value = 0.021110617731369206
$ws.Range("A4").Value = 0.020401408548431157
$ws.Range("A5").Value = 0.020154365035892349
$ws.Range("A6").Value = 0.018460373398173722
